# Bugfixed the naive forecaster component module
#
# The naive AR(2) forecaster's first observation row (the stub row that
# only carried the base date/year with no forecast yet) was dropped, and
# the y_1_forecast (column E) values for the remaining rows were
# recomputed by the fixed forecaster logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete leading data row (old row 2); everything below
# shifts up by one, which is exactly what the diff shows for columns
# A-D (dates/years/y_0_forecast all shift up one row, values unchanged).
$ws.Rows.Item(2).Delete()

# Recomputed y_1_forecast (column E) values. Rows 2-5 (post-shift) have
# no forecast yet, so clear any leftover values there; rows 6-18 get the
# corrected forecast figures.
$ws.Range("E2:E5").ClearContents()

$e = @{
    6  = 0.4163953164477929
    7  = 1.653207170606596
    8  = 1.270027657109818
    9  = 1.579162878845075
    10 = 1.452243308058287
    11 = 2.068578555939404
    12 = 1.651658474923545
    13 = 0.115841687688345
    14 = -2.092304328310923
    15 = 1.533339625605379
    16 = 0.492911192428136
    17 = 0.2100922168233987
    18 = 0.5208382580577098
}

foreach ($row in $e.Keys) {
    $ws.Cells.Item($row, 5).Value = $e[$row]
}

# Fix up the tiny floating-point nudge on C3 (y_0_forecast) introduced by
# the recompute.
$ws.Cells.Item(3, 3).Value = -1.611885206309638
